$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet's tab/title text (the workbook.xml <sheet name="..."> entry)
$ws.Name = "Through 2022-08-22"

# Update the "August" row label in column A (shared string)
$ws.Range("A9").Value = "August (through 08-22)"

# Update August row (row 9) values for years 2015-2022 (columns B-I)
$ws.Range("B9").Value = 25
$ws.Range("C9").Value = 51
$ws.Range("D9").Value = 59
$ws.Range("E9").Value = 39
$ws.Range("F9").Value = 29
$ws.Range("G9").Value = 133
$ws.Range("H9").Value = 115
$ws.Range("I9").Value = 129

# Update Total row (row 10) values for years 2015-2022 (columns B-I)
$ws.Range("B10").Value = 187
$ws.Range("C10").Value = 353
$ws.Range("D10").Value = 524
$ws.Range("E10").Value = 464
$ws.Range("F10").Value = 333
$ws.Range("G10").Value = 754
$ws.Range("H10").Value = 1025
$ws.Range("I10").Value = 1100
